# Auto-update Data Telemetría - Ejecución Diaria
# Appends the newest daily connectivity snapshot to the "historico" sheet
# and refreshes the "ultimo_snapshot" sheet with that same latest data.

$wb = $excel.ActiveWorkbook

$historico = $wb.Worksheets.Item("historico")
$ultimoSnapshot = $wb.Worksheets.Item("ultimo_snapshot")

# New rows of data for the latest run date.
$fecha = 45992

$nuevaData = @(
    @($fecha, "Telemetría",                     5902, 3459, 590, 178, 673, 1002, 58.61, 10,   3.02, 11.4,  16.98),
    @($fecha, "GPS (según REGLA)",               5302, 4688, 344,  88, 176,    6, 88.42, 6.49, 1.66,  3.32,  0.11),
    @($fecha, "GPS (todas con gps_timestamp)",  11198, 9476, 862, 279, 581,    0, 84.62, 7.7,  2.49,  5.19,  0)
)

# Append the new rows to the bottom of "historico".
$startRow = $historico.Cells.Item($historico.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $nuevaData.Count; $i++) {
    $row = $startRow + $i
    $values = $nuevaData[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $historico.Cells.Item($row, $c + 1).Value = $values[$c]
    }
    # Match the "fecha" column date formatting used by the existing rows.
    $historico.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
}

# Overwrite "ultimo_snapshot" rows 2-4 with this same latest data.
for ($i = 0; $i -lt $nuevaData.Count; $i++) {
    $row = 2 + $i
    $values = $nuevaData[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ultimoSnapshot.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
